$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCommitment")

$ws.Range("K1").Value = "From Currency"
$ws.Range("L1").Value = "To Currency"
$ws.Range("M1").Value = "Exchange Rate "
$ws.Range("N1").Value = "As Of"

$ws.Range("K2").Select()
